# Applies the commit's changes:
#  1. Updates the "Förändrad" date in column C (for every data row) from
#     45184 to 45186.
#  2. Adds the friendly display-text argument (the Beteckning value from
#     column A of the same row) as the second HYPERLINK() argument for the
#     link-formula columns S, T, V, W, X, Y wherever such a formula exists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row (header is row 1, data starts row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Columns that may hold a HYPERLINK() formula needing the new 2nd argument.
$linkCols = @(19, 20, 22, 23, 24, 25)   # S, T, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {

    # --- 1. Bump the "Förändrad" date in column C ---------------------
    $cCell = $ws.Cells.Item($r, 3)
    $cVal = $cCell.Value2()
    if ($cVal -ne $null -and [double]$cVal -eq 45184) {
        $cCell.Value = 45186
    }

    # --- 2. Add display text to HYPERLINK formulas in this row --------
    $aVal = $ws.Cells.Item($r, 1).Value()

    foreach ($colIdx in $linkCols) {
        $cell = $ws.Cells.Item($r, $colIdx)
        if ($cell.HasFormula) {
            $f = $cell.Formula()
            if ($f -and $f.Length -gt 0 -and $f.Substring($f.Length - 1) -eq ")") {
                $newF = $f.Substring(0, $f.Length - 1) + ', "' + $aVal + '")'
                $cell.Formula = $newF
            }
        }
    }
}
